$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")
$ws.Activate()

# --- Data edits: final exam score (O) for two students, formulas recalc automatically ---
$ws.Range("O10").Value = 19
$ws.Range("O13").Value = 28

# --- Column Q (17) width -> 5 (best-fit-like width) ---
$ws.Columns("Q").ColumnWidth = 4.14

# --- Page setup: print scale 82 -> 87 ---
$ps = $ws.PageSetup
$ps.Zoom = 87
$ps.FitToPagesTall = $false
$ps.FitToPagesWide = 1

# --- View: scroll position and active selection ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Y9").Select()
